$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-21 Sunday" "2024-01-22 Monday"

Replace-Text "145×8=1160" "367×8=2936"
Replace-Text "638×9=5742" "854×8=6832"
Replace-Text "800×8=6400" "664×8=5312"
Replace-Text "612×9=5508" "734×8=5872"
Replace-Text "267×5=1335" "601×5=3005"
Replace-Text "713×8=5704" "818×7=5726"
Replace-Text "375×8=3000" "233×3=699"
Replace-Text "874×2=1748" "546×2=1092"
Replace-Text "188×3=564" "967×6=5802"
Replace-Text "458×6=2748" "735×4=2940"
Replace-Text "129×7=903" "562×2=1124"
Replace-Text "341×6=2046" "142×5=710"
Replace-Text "257×6=1542" "636×6=3816"
Replace-Text "142×6=852" "667×3=2001"
Replace-Text "566×4=2264" "619×5=3095"
Replace-Text "904×4=3616" "534×2=1068"
Replace-Text "744×5=3720" "521×9=4689"
Replace-Text "750×7=5250" "886×9=7974"
Replace-Text "450×9=4050" "840×7=5880"
Replace-Text "498×7=3486" "118×2=236"
Replace-Text "723×7=5061" "305×8=2440"
Replace-Text "978×9=8802" "886×6=5316"
Replace-Text "417×7=2919" "102×5=510"
Replace-Text "984×5=4920" "341×5=1705"
Replace-Text "412×9=3708" "357×6=2142"
